$wb = $excel.ActiveWorkbook

# ---- Step 1: turn the existing '总计' sheet into '2022-Q1' (reuses sheetId/rId) ----
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# ---- Step 2: new header row (basi代码/名称/规模/仓位/占比/市值/排名) ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---- Step 3: extend the A-column index style (s=2) down to row 25 ----
$q1.Range("A2").Copy()
$q1.Range("A3:A25").PasteSpecial(-4122)

# ---- Step 4: force D:G to be entered as text (matches source inlineStr data), like the original author typing it in as text ----
$q1.Range("D2:G25").NumberFormat = "@"
$q1.Range("B2:B25").NumberFormat = "@"

# ---- Step 5: fill in the 24 fund rows ----
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "001838"
$q1.Range("C2").Value = "国投瑞银国家安全灵活配置混合"
$q1.Range("D2").Value = "32.13"
$q1.Range("E2").Value = "94.68"
$q1.Range("F2").Value = "7.80"
$q1.Range("G2").Value = "2.5061"
$q1.Range("H2").Value = 5

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "519196"
$q1.Range("C3").Value = "万家新兴蓝筹灵活配置混合"
$q1.Range("D3").Value = "21.26"
$q1.Range("E3").Value = "80.70"
$q1.Range("F3").Value = "3.26"
$q1.Range("G3").Value = "0.6931"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "009199"
$q1.Range("C4").Value = "万家价值优势一年持有期混合"
$q1.Range("D4").Value = "13.70"
$q1.Range("E4").Value = "89.53"
$q1.Range("F4").Value = "3.52"
$q1.Range("G4").Value = "0.4822"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "340001"
$q1.Range("C5").Value = "兴全可转债混合"
$q1.Range("D5").Value = "43.74"
$q1.Range("E5").Value = "27.50"
$q1.Range("F5").Value = "1.04"
$q1.Range("G5").Value = "0.4549"
$q1.Range("H5").Value = 3

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "630001"
$q1.Range("C6").Value = "华商领先企业混合"
$q1.Range("D6").Value = "9.21"
$q1.Range("E6").Value = "83.85"
$q1.Range("F6").Value = "3.69"
$q1.Range("G6").Value = "0.3398"
$q1.Range("H6").Value = 3

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "000690"
$q1.Range("C7").Value = "前海开源大海洋战略经济灵活配置混合"
$q1.Range("D7").Value = "6.91"
$q1.Range("E7").Value = "93.87"
$q1.Range("F7").Value = "4.16"
$q1.Range("G7").Value = "0.2875"
$q1.Range("H7").Value = 8

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "168501"
$q1.Range("C8").Value = "北信瑞丰产业升级多策略混合"
$q1.Range("D8").Value = "4.42"
$q1.Range("E8").Value = "94.11"
$q1.Range("F8").Value = "5.25"
$q1.Range("G8").Value = "0.2320"
$q1.Range("H8").Value = 3

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "008009"
$q1.Range("C9").Value = "华商高端装备制造股票"
$q1.Range("D9").Value = "5.40"
$q1.Range("E9").Value = "88.04"
$q1.Range("F9").Value = "3.96"
$q1.Range("G9").Value = "0.2138"
$q1.Range("H9").Value = 3

$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "460002"
$q1.Range("C10").Value = "华泰柏瑞积极成长混合A"
$q1.Range("D10").Value = "6.11"
$q1.Range("E10").Value = "81.55"
$q1.Range("F10").Value = "3.20"
$q1.Range("G10").Value = "0.1955"
$q1.Range("H10").Value = 8

$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "010460"
$q1.Range("C11").Value = "兴业研究精选混合"
$q1.Range("D11").Value = "3.41"
$q1.Range("E11").Value = "89.54"
$q1.Range("F11").Value = "2.82"
$q1.Range("G11").Value = "0.0962"
$q1.Range("H11").Value = 10

$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "010375"
$q1.Range("C12").Value = "国金鑫悦经济新动能混合A"
$q1.Range("D12").Value = "1.29"
$q1.Range("E12").Value = "94.93"
$q1.Range("F12").Value = "7.18"
$q1.Range("G12").Value = "0.0926"
$q1.Range("H12").Value = 4

$q1.Range("A13").Value = 11
$q1.Range("B13").Value = "001457"
$q1.Range("C13").Value = "华商新常态灵活配置混合"
$q1.Range("D13").Value = "2.32"
$q1.Range("E13").Value = "83.98"
$q1.Range("F13").Value = "3.69"
$q1.Range("G13").Value = "0.0856"
$q1.Range("H13").Value = 3

$q1.Range("A14").Value = 12
$q1.Range("B14").Value = "007923"
$q1.Range("C14").Value = "方正富邦天鑫灵活配置混合A"
$q1.Range("D14").Value = "1.14"
$q1.Range("E14").Value = "80.49"
$q1.Range("F14").Value = "5.29"
$q1.Range("G14").Value = "0.0603"
$q1.Range("H14").Value = 4

$q1.Range("A15").Value = 13
$q1.Range("B15").Value = "007924"
$q1.Range("C15").Value = "方正富邦天鑫灵活配置混合C"
$q1.Range("D15").Value = "1.14"
$q1.Range("E15").Value = "80.49"
$q1.Range("F15").Value = "5.29"
$q1.Range("G15").Value = "0.0603"
$q1.Range("H15").Value = 4

$q1.Range("A16").Value = 14
$q1.Range("B16").Value = "673050"
$q1.Range("C16").Value = "西部利得新盈灵活配置混合"
$q1.Range("D16").Value = "2.28"
$q1.Range("E16").Value = "68.90"
$q1.Range("F16").Value = "1.95"
$q1.Range("G16").Value = "0.0445"
$q1.Range("H16").Value = 8

$q1.Range("A17").Value = 15
$q1.Range("B17").Value = "002604"
$q1.Range("C17").Value = "华夏新起点灵活配置混合A"
$q1.Range("D17").Value = "0.90"
$q1.Range("E17").Value = "75.30"
$q1.Range("F17").Value = "3.39"
$q1.Range("G17").Value = "0.0305"
$q1.Range("H17").Value = 6

$q1.Range("A18").Value = 16
$q1.Range("B18").Value = "970015"
$q1.Range("C18").Value = "申万宏源红利成长灵活配置混合"
$q1.Range("D18").Value = "1.04"
$q1.Range("E18").Value = "70.39"
$q1.Range("F18").Value = "2.46"
$q1.Range("G18").Value = "0.0256"
$q1.Range("H18").Value = 7

$q1.Range("A19").Value = 17
$q1.Range("B19").Value = "010376"
$q1.Range("C19").Value = "国金鑫悦经济新动能混合C"
$q1.Range("D19").Value = "0.28"
$q1.Range("E19").Value = "94.93"
$q1.Range("F19").Value = "7.18"
$q1.Range("G19").Value = "0.0201"
$q1.Range("H19").Value = 4

$q1.Range("A20").Value = 18
$q1.Range("B20").Value = "004926"
$q1.Range("C20").Value = "中航军民融合精选混合A"
$q1.Range("D20").Value = "0.35"
$q1.Range("E20").Value = "91.27"
$q1.Range("F20").Value = "4.80"
$q1.Range("G20").Value = "0.0168"
$q1.Range("H20").Value = 9

$q1.Range("A21").Value = 19
$q1.Range("B21").Value = "004927"
$q1.Range("C21").Value = "中航军民融合精选混合C"
$q1.Range("D21").Value = "0.27"
$q1.Range("E21").Value = "91.27"
$q1.Range("F21").Value = "4.80"
$q1.Range("G21").Value = "0.0130"
$q1.Range("H21").Value = 9

$q1.Range("A22").Value = 20
$q1.Range("B22").Value = "010663"
$q1.Range("C22").Value = "长江均衡成长混合A"
$q1.Range("D22").Value = "0.26"
$q1.Range("E22").Value = "85.90"
$q1.Range("F22").Value = "4.85"
$q1.Range("G22").Value = "0.0126"
$q1.Range("H22").Value = 2

$q1.Range("A23").Value = 21
$q1.Range("B23").Value = "010664"
$q1.Range("C23").Value = "长江均衡成长混合C"
$q1.Range("D23").Value = "0.05"
$q1.Range("E23").Value = "85.90"
$q1.Range("F23").Value = "4.85"
$q1.Range("G23").Value = "0.0024"
$q1.Range("H23").Value = 2

$q1.Range("A24").Value = 22
$q1.Range("B24").Value = "008213"
$q1.Range("C24").Value = "华夏新起点灵活配置混合C"
$q1.Range("D24").Value = "0.02"
$q1.Range("E24").Value = "75.30"
$q1.Range("F24").Value = "3.39"
$q1.Range("G24").Value = "0.0007"
$q1.Range("H24").Value = 6

$q1.Range("A25").Value = 23
$q1.Range("B25").Value = "960030"
$q1.Range("C25").Value = "华泰柏瑞积极成长混合H"
$q1.Range("D25").Value = "0.00"
$q1.Range("E25").Value = "81.55"
$q1.Range("F25").Value = "3.20"
$q1.Range("G25").NumberFormat = "General"
$q1.Range("G25").Value = 0
$q1.Range("H25").Value = 8

# ---- Step 6: add the new '总计' sheet right after '2022-Q1' (gets next sheetId/rId) ----
$after = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Add($null, $after)
$total.Name = "总计"

# ---- Step 7: recreate the historical-totals table, with the new 2022-Q1 row on top ----
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 24
$total.Range("D2").Value = 5.97

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 19
$total.Range("D3").Value = 6.79

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 13
$total.Range("D4").Value = 7.59

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 12
$total.Range("D5").Value = 6.83

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.06

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 12
$total.Range("D7").Value = 0.26

# ---- Step 8: restore the original active sheet/tab (unrelated to this edit) ----
$wb.Worksheets.Item("2020-Q4").Activate()

